{"js": "// The document has two near-identical \"Restart the stage...\" sections\n// (one already highlighted darkGray under numId=3's sibling list, and the\n// target one - under the numId=3 list flanked by cyan-highlighted items -\n// which currently has no highlight). We need to:\n//   1. Highlight (yellow) the \"Restart the stage when Robetroid dies...\"\n//      paragraph, the \"Depleted health\" paragraph, and the \"Falls \"\n//      paragraph (paragraph mark + every run).\n//   2. Move the \"_GoBack\" bookmark from the end of the\n//      \"... when touching him\" paragraph to the end of the \"Falls \"\n//      paragraph.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Load highlight color for every paragraph up front so the search below can\n// pick the correct (un-highlighted) triplet in a single pass.\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].font.load(\"highlightColor\");\n}\nawait context.sync();\n\n// Locate the \"Restart the stage... / Depleted health / Falls \" triplet that\n// is NOT already highlighted (the darkGray copy earlier in the doc - same\n// text, already highlighted - must stay untouched).\nlet restartPara = null;\nlet depletedPara = null;\nlet fallsPara = null;\n\nfor (let i = 0; i + 2 < paragraphs.items.length; i++) {\n  const p0 = paragraphs.items[i];\n  const p1 = paragraphs.items[i + 1];\n  const p2 = paragraphs.items[i + 2];\n  if (\n    p0.text === \"Restart the stage when Robetroid dies as long as there are extra lives\" &&\n    p1.text === \"Depleted health\" &&\n    p2.text === \"Falls \" &&\n    !p0.font.highlightColor\n  ) {\n    restartPara = p0;\n    depletedPara = p1;\n    fallsPara = p2;\n    break;\n  }\n}\n\nif (!restartPara) {\n  throw new Error(\"Could not locate the un-highlighted Restart/Depleted/Falls triplet\");\n}\n\n// 1. Apply yellow highlight to the paragraph mark + every run in each of\n// the three paragraphs (setting paragraph.font.highlightColor covers the\n// whole paragraph range, including the trailing paragraph mark run).\nrestartPara.font.highlightColor = \"Yellow\";\ndepletedPara.font.highlightColor = \"Yellow\";\nfallsPara.font.highlightColor = \"Yellow\";\nawait context.sync();\n\n// 2. Move the \"_GoBack\" bookmark to the end of the \"Falls \" paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst fallsEnd = fallsPara.getRange(\"End\");\nfallsEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The document contains two near-identical \"Restart the stage...\" blocks\n# (one already highlighted darkGray, and the target block - flanked by\n# cyan-highlighted list items - which currently has no highlight at all).\n# We need to:\n#   1. Highlight (yellow) the \"Restart the stage when Robetroid dies...\"\n#      paragraph, the \"Depleted health\" paragraph, and the \"Falls \"\n#      paragraph (this covers both the paragraph mark / pilcrow and every\n#      run in each paragraph).\n#   2. Move the \"_GoBack\" bookmark from the end of the\n#      \"... when touching him\" paragraph to the end of the \"Falls \"\n#      paragraph (right after its text, before its paragraph mark).\n\n$d = $word.ActiveDocument\n\n# --- Locate the un-highlighted \"Restart the stage / Depleted health / Falls\"\n# --- triplet (skip the darkGray copy earlier in the document). A Paragraph's\n# --- Range.Text always carries its trailing paragraph-mark (\"`r\"), so we\n# --- compare against that form.\n$restartText  = \"Restart the stage when Robetroid dies as long as there are extra lives`r\"\n$depletedText = \"Depleted health`r\"\n$fallsText    = \"Falls `r\"\n\n$total = $d.Paragraphs.Count\n$restartIdx = -1\nfor ($i = 1; $i -le ($total - 2); $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -eq $restartText) {\n        $hl = $d.Paragraphs.Item($i).Range.HighlightColorIndex\n        $next1 = $d.Paragraphs.Item($i + 1).Range.Text\n        $next2 = $d.Paragraphs.Item($i + 2).Range.Text\n        if ((-not $hl) -and $next1 -eq $depletedText -and $next2 -eq $fallsText) {\n            $restartIdx = $i\n            break\n        }\n    }\n}\n\nif ($restartIdx -eq -1) {\n    throw \"Could not locate the un-highlighted Restart/Depleted/Falls paragraph triplet\"\n}\n\n$restartPara  = $d.Paragraphs.Item($restartIdx)\n$depletedPara = $d.Paragraphs.Item($restartIdx + 1)\n$fallsPara    = $d.Paragraphs.Item($restartIdx + 2)\n\n# --- 1. Apply yellow highlight (wdYellow = 7) to each paragraph's full\n# --- range, which covers the paragraph mark as well as every run.\n$restartPara.Range.Font.HighlightColorIndex  = 7\n$depletedPara.Range.Font.HighlightColorIndex = 7\n$fallsPara.Range.Font.HighlightColorIndex    = 7\n\n# --- 2. Move the \"_GoBack\" bookmark to the end of the \"Falls \" paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Re-fetch the Falls paragraph (indices are unaffected by the highlight\n# edits above, but this keeps the reference fresh/explicit).\n$fallsPara = $d.Paragraphs.Item($restartIdx + 2)\n\n# Paragraph.Range.Text includes the trailing paragraph mark, so its length\n# minus 1 is the offset of the position right after all visible text and\n# right before that mark - exactly where the bookmark needs to land.\n$endOfTextPos = $fallsPara.Range.Start + $fallsPara.Range.Text.Length - 1\n\n# Bookmarking a position that sits immediately before a paragraph mark can\n# collide with the mark itself, so anchor a scratch character there first,\n# insert the bookmark right before the scratch character, then remove the\n# scratch character again.\n$scratchMarker = [char]1\n$insertPoint = $d.Range($endOfTextPos, $endOfTextPos)\n$insertPoint.InsertAfter($scratchMarker)\n\n$bookmarkRange = $d.Range($endOfTextPos, $endOfTextPos)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n\n$scratchRange = $d.Range($endOfTextPos, $endOfTextPos + 1)\n$scratchRange.Delete()\n"}
